$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text (translation) corrections in column D ---

$ws.Range("D4").Value = " అత్యధిక వికెట్లు వికెట్ కీపర్ క్యాచ్ ల ద్వారా సాధించిన ఆటగాళ్ల"
$ws.Range("D7").Value = " ఒకే ఇన్నింగ్స్ లో ఫోర్లు, సిక్సర్ల ద్వారా అత్యధిక పరుగులు సాధించిన ఆటగాళ్ల "
$ws.Range("D10").Value = " కెరీర్ లో ఒక ఇన్నింగ్స్ లో నాలుగు వికెట్లు అత్యధిక సార్లు సాధించిన ఆటగాళ్ల "
$ws.Range("D11").Value = " అతి తక్కువ కాలం నివసించిన ఆటగాళ్ల"
$ws.Range("D13").Value = "అదనపు పరుగులు లేకుండా అత్యధిక ఇన్నింగ్స్ మొత్తం సమర్పించిన వికెట్ కీపర్ల"
$ws.Range("D16").Value = " వరుస ఇన్నింగ్స్ లలో నాలుగు వికెట్లు సాధించిన ఆటగాళ్ల"
$ws.Range("D19").Value = " ఒక జట్టుకు ఆడిన రెండు ప్రదర్శనల మధ్య అత్యధిక వరుస మ్యాచ్ లు ఆడకుండా ఉన్న ఆటగాళ్ల "
$ws.Range("D31").Value = " రెండో వికెట్ కు అత్యధిక భాగస్వామ్యం వహించిన ఆటగాళ్ల"
$ws.Range("D39").Value = " సజీవంగా ఉన్న పురాతన ( ఓల్డ్ ) ఆటగాళ్ల"
$ws.Range("D42").Value = " కెప్టెన్ గా అరంగేట్రం చేసిన పురాతన ( ఓల్డ్ ) కెప్టెన్ల"
$ws.Range("D50").Value = "ఒక మ్యాచ్ లో ఓడిపోయిన జట్టుకి అత్యధిక పరుగులు చేసిన ఆటగాళ్ల"
$ws.Range("D59").Value = " శతకం సాధించిన పురాతన ( ఓల్డ్ ) ఆటగాళ్ల "
$ws.Range("D61").Value = " ఒక ఇన్నింగ్స్ లో అత్యధిక అదనపు పరుగులు(బైస్) సమర్పించిన ఆటగాళ్ల  "
$ws.Range("D66").Value = " ఆరో వికెట్ కు అత్యధిక భాగస్వామ్యం చేసిన ఆటగాళ్ల"
$ws.Range("D70").Value = " ఉత్తమ కెరీర్ బౌలింగ్ సగటు ఉన్న ఆటగాళ్ల"
$ws.Range("D71").Value = " ఒక మ్యాచ్ లో అత్యధిక క్యాచ్ లు పట్టిన ఆటగాళ్ల  "

# --- Formatting: align D11, D31, D32, D66, D70, D72 with the same style
#     already used by similar translation cells (e.g. D3) ---

$ws.Range("D3").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("D66").PasteSpecial(-4122)
$ws.Range("D70").PasteSpecial(-4122)
$ws.Range("D72").PasteSpecial(-4122)
